# Update "想去人数" (Column F) values on the "展览" and "全部类型" sheets
# to reflect newly scraped counts, per the commit:
#   "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) receive the same updates,
# except row 22 which differs by one between the two sheets.
$wsZhanLan = $wb.Worksheets.Item("展览")
$wsQuanBu  = $wb.Worksheets.Item("全部类型")

# Row -> new F-column value, identical for both sheets except row 22.
$common = @{
    2  = 277
    3  = 1398
    4  = 166
    6  = 238
    11 = 4692
    12 = 6948
    14 = 60
    16 = 577
    18 = 4164
    19 = 871
    21 = 67
    26 = 381
    27 = 379
    29 = 240
    30 = 48
    31 = 1647
    32 = 1047
    33 = 70
    34 = 489
    37 = 4
    41 = 206
    42 = 650
}

foreach ($ws in @($wsZhanLan, $wsQuanBu)) {
    foreach ($row in $common.Keys) {
        $ws.Cells.Item($row, 6).Value = $common[$row]
    }
}

# Row 22 differs between the two sheets.
$wsZhanLan.Cells.Item(22, 6).Value = 2744
$wsQuanBu.Cells.Item(22, 6).Value = 2745
